{"js": "// The document contains a few \"Ripete dal passo N finch\u00e9 ... soddisfatt[oa]\"\n// paragraphs. Some of them are already a single run; three of them are split\n// across several runs (separated by <w:proofErr> grammar-check markers) as\n// a left-over from Word's grammar checker, e.g.:\n//   run1: \"Ripete dal passo \"\n//   proofErr gramStart\n//   run2: \"6\"\n//   proofErr gramEnd\n//   run3: \" finch\u00e9 non \u00e8 soddisfatto\"\n// The edit collapses each such paragraph back into a single run (and drops\n// the now meaningless proofErr markers) without changing the visible text\n// or its (italic) formatting.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Paragraphs that still carry the split runs use the accented \"\u00e9\" in\n// \"finch\u00e9\" (the ones already collapsed into a single run use the plain\n// grave accent \"finch\u00e8\" instead) - use that, together with the leading\n// wording, to find exactly the paragraphs that must be merged.\nconst targetPattern = /^Ripete dal passo \\d+ finch\u00e9 non \u00e8 soddisfatt[oa]$/;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text;\n\n  if (!targetPattern.test(text)) {\n    continue;\n  }\n\n  // Normalize any internal run-splitting whitespace quirks and rewrite the\n  // whole paragraph range as one run; Word.js merges the replacement into a\n  // single run using the formatting already present at that location\n  // (italic is preserved) and removes the old proofErr markers along with\n  // the runs they used to separate.\n  const range = paragraph.getRange();\n  range.insertText(text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document contains a few \"Ripete dal passo N finch\u00e9 ... soddisfatt[oa]\"\n# paragraphs. Some of them already live in a single run; three of them are\n# still split across several runs (separated by grammar-check <w:proofErr>\n# markers) as a left-over from Word's grammar checker, e.g.:\n#   run1: \"Ripete dal passo \"\n#   proofErr gramStart\n#   run2: \"6\"\n#   proofErr gramEnd\n#   run3: \" finch\u00e9 non \u00e8 soddisfatto\"\n# This collapses each such paragraph back into a single run (dropping the\n# now meaningless proofErr markers) without changing the visible text or\n# its (italic) formatting.\n\n$d = $word.ActiveDocument\n\n# Paragraphs that still carry the split runs use the accented \"\u00e9\" in\n# \"finch\u00e9\" (the ones already collapsed into a single run use the plain\n# grave accent \"finch\u00e8\" instead) - together with the leading wording this\n# reliably singles out exactly the paragraphs that must be merged. The\n# paragraph text can end either with a plain paragraph mark (Cr) or, for\n# the last paragraph of a table cell, a Cr followed by a cell mark (Bel).\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -match \"^Ripete dal passo \\d+ finch\u00e9 non \u00e8 soddisfatt[oa]\\r[\\a]?$\") {\n        $target = $t.TrimEnd([char]13, [char]7)\n\n        $r = $p.Range\n        $r.Find.ClearFormatting()\n        $r.Find.Replacement.ClearFormatting()\n        # Replacing the exact text in place merges it back into a single\n        # run (using the formatting already present there, so italics are\n        # kept) and removes the proofErr markers that used to split it.\n        # Wrap:=wdFindStop (0) and Replace:=wdReplaceOne (1) keep the\n        # operation confined to this paragraph's own range instead of\n        # touching a later paragraph that happens to contain the same text.\n        $r.Find.Execute($target, $false, $false, $false, $false, $false, $true, 0, $false, $target, 1) | Out-Null\n    }\n}\n"}
